$d = $word.ActiveDocument

# Locate the unique paragraph "What are the sub-goals?" that belongs to the
# "Socks in the Dark" problem -- i.e. the one immediately following the
# paragraph that ends with "Trying to match socks in the dark."
$target = $null
$found = $false
foreach ($p in $d.Paragraphs) {
    if ($found) {
        $target = $p
        break
    }
    if ($p.Range.Text -like "*Trying to match socks in the dark.*") {
        $found = $true
    }
}

if ($target -eq $null) {
    throw "Could not locate the 'What are the sub-goals?' paragraph for the Socks in the Dark problem."
}

$r = $target.Range

# Append " Trying to match three pairs, one of each color." at the end of
# the paragraph (i.e. right before the terminating paragraph mark).
$r.InsertAfter(" ")
$coloredStart = $r.End - 1
$r.InsertAfter("Trying to match three pairs, one of each color.")
$coloredEnd = $r.End - 1

$coloredRange = $d.Range($coloredStart, $coloredEnd)
$coloredRange.Font.Color = 0xFF6633   # wdColor BGR value for RGB(0x33,0x66,0xFF)

# The "_GoBack" bookmark tracks the location of the most recent edit, so it
# moves from around the old "Trying to match socks in the dark." text to
# wrap the newly-typed "Trying to match three pairs, one of each color." text.
$old = $d.Bookmarks("_GoBack")
$old.Delete()
$d.Bookmarks.Add("_GoBack", $coloredRange)
